# Auto-generated Excel COM script applying value updates per the commit diff
# (scheduled runner update to Masamune_Profits.xlsx Leve Profit sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink
$ws.Range("H28").Value = 471.3889
$ws.Range("I28").Value = 384.91666
$ws.Range("J28").Value = 644.3333
$ws.Range("K28").Value = 384.91666
$ws.Range("L28").Value = 644.3333
$ws.Range("M28").Value = 100.08334
$ws.Range("N28").Value = -1614.3333

# Row 80: Cleansing the Wicked Humours | Hallowed Water
$ws.Range("H80").Value = 10417167
$ws.Range("I80").Value = 18519080
$ws.Range("J80").Value = 422.2143
$ws.Range("K80").Value = 55557240
$ws.Range("L80").Value = 1266.6429
$ws.Range("M80").Value = -55556242
$ws.Range("N80").Value = -3262.6429

# Row 83: Washing Away the Sins (L) | Hallowed Water
$ws.Range("H83").Value = 10417167
$ws.Range("I83").Value = 18519080
$ws.Range("J83").Value = 422.2143
$ws.Range("K83").Value = 166671720
$ws.Range("L83").Value = 3799.9287
$ws.Range("M83").Value = -166666728
$ws.Range("N83").Value = -13783.9287

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 23585.596
$ws.Range("I132").Value = 3031.3713
$ws.Range("J132").Value = 126356.71
$ws.Range("K132").Value = 9094.1139
$ws.Range("L132").Value = 379070.13
$ws.Range("M132").Value = -6564.1139
$ws.Range("N132").Value = -384130.13

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2335309.2
$ws.Range("I137").Value = 7693263
$ws.Range("J137").Value = 5764.1304
$ws.Range("K137").Value = 23079789
$ws.Range("L137").Value = 17292.3912
$ws.Range("M137").Value = -23077239
$ws.Range("N137").Value = -22392.3912

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 2672.25
$ws.Range("J141").Value = 5737.5
$ws.Range("L141").Value = 17212.5
$ws.Range("N141").Value = -27572.5

$ws = $wb.Worksheets.Item("ARM")
# Row 29: No Hand-me-downs | Iron Vambraces
$ws.Range("H29").Value = 6612.5
$ws.Range("I29").Value = 966.6667
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 966.6667
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -658.6667
$ws.Range("N29").Value = -10616

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 8977.145
$ws.Range("I32").Value = 7657.712
$ws.Range("J32").Value = 38004.668
$ws.Range("K32").Value = 7657.712
$ws.Range("L32").Value = 38004.668
$ws.Range("M32").Value = -7370.712
$ws.Range("N32").Value = -38578.668

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1846.0476
$ws.Range("I45").Value = 1818.4546
$ws.Range("J45").Value = 1876.4
$ws.Range("K45").Value = 1818.4546
$ws.Range("L45").Value = 1876.4
$ws.Range("M45").Value = -1441.4546
$ws.Range("N45").Value = -2630.4

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1787.4565
$ws.Range("I74").Value = 1441.0571
$ws.Range("J74").Value = 2889.6365
$ws.Range("K74").Value = 1441.0571
$ws.Range("L74").Value = 2889.6365
$ws.Range("M74").Value = -567.0571
$ws.Range("N74").Value = -4637.636500000001

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1787.4565
$ws.Range("I77").Value = 1441.0571
$ws.Range("J77").Value = 2889.6365
$ws.Range("K77").Value = 7205.2855
$ws.Range("L77").Value = 14448.1825
$ws.Range("M77").Value = -2837.2855
$ws.Range("N77").Value = -23184.1825

# Row 120: One Foot Forward | Dwarven Mythril Shoes of Maiming
$ws.Range("H120").Value = 43420
$ws.Range("J120").Value = 43420
$ws.Range("L120").Value = 43420
$ws.Range("N120").Value = -53096

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2425.25
$ws.Range("I132").Value = 1346.1818
$ws.Range("K132").Value = 4038.5454
$ws.Range("M132").Value = -1508.5454

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2236.3333
$ws.Range("I134").Value = 1394.12
$ws.Range("K134").Value = 4182.36
$ws.Range("M134").Value = -1647.36

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 4503.09
$ws.Range("I31").Value = 1921.16
$ws.Range("K31").Value = 1921.16
$ws.Range("M31").Value = -1626.16

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 4503.09
$ws.Range("I34").Value = 1921.16
$ws.Range("K34").Value = 1921.16
$ws.Range("M34").Value = -1719.16

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1690.5476
$ws.Range("I58").Value = 1111.8518
$ws.Range("J58").Value = 2732.2
$ws.Range("K58").Value = 1111.8518
$ws.Range("L58").Value = 2732.2
$ws.Range("M58").Value = -908.8517999999999
$ws.Range("N58").Value = -3138.2

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1690.5476
$ws.Range("I136").Value = 1111.8518
$ws.Range("J136").Value = 2732.2
$ws.Range("K136").Value = 3335.5554
$ws.Range("L136").Value = 8196.599999999999
$ws.Range("M136").Value = -785.5553999999997
$ws.Range("N136").Value = -13296.6

$ws = $wb.Worksheets.Item("CUL")
# Row 17: Chew the Fat | Grilled Dodo
$ws.Range("H17").Value = 800
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 975
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 2925
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -3263

# Row 40: True Grits | Cornmeal
$ws.Range("H40").Value = 4956.273
$ws.Range("I40").Value = 6375.875
$ws.Range("J40").Value = 1170.6666
$ws.Range("K40").Value = 25503.5
$ws.Range("L40").Value = 4682.6664
$ws.Range("M40").Value = -25434.5
$ws.Range("N40").Value = -4820.6664

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 877.85
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 877.85
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2633.55
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = -12713.55

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2817.8
$ws.Range("J132").Value = 3704.4
$ws.Range("L132").Value = 33339.6
$ws.Range("N132").Value = -38399.6

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 4799.0586
$ws.Range("I70").Value = 4821.893
$ws.Range("J70").Value = 4692.5
$ws.Range("K70").Value = 4821.893
$ws.Range("L70").Value = 4692.5
$ws.Range("M70").Value = -4551.893
$ws.Range("N70").Value = -5232.5

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 4799.0586
$ws.Range("I73").Value = 4821.893
$ws.Range("J73").Value = 4692.5
$ws.Range("K73").Value = 4821.893
$ws.Range("L73").Value = 4692.5
$ws.Range("M73").Value = -3885.893
$ws.Range("N73").Value = -6564.5

# Row 110: Slimming Down | Stonegold Rapier
$ws.Range("H110").Value = 41997.668
$ws.Range("J110").Value = 41997.668
$ws.Range("L110").Value = 41997.668
$ws.Range("N110").Value = -50177.668

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 38466508
$ws.Range("I132").Value = 76928770
$ws.Range("J132").Value = 4246.615
$ws.Range("K132").Value = 230786310
$ws.Range("L132").Value = 12739.845
$ws.Range("M132").Value = -230783780
$ws.Range("N132").Value = -17799.845

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 3084.5715
$ws.Range("I7").Value = 2116.125
$ws.Range("J7").Value = 4375.8335
$ws.Range("K7").Value = 2116.125
$ws.Range("L7").Value = 4375.8335
$ws.Range("M7").Value = -2004.125
$ws.Range("N7").Value = -4599.8335

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 333
$ws.Range("I22").Value = 351.33334
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 351.33334
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -56.33334000000002
$ws.Range("N22").Value = -890

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 333
$ws.Range("I27").Value = 351.33334
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 351.33334
$ws.Range("L27").Value = 300
$ws.Range("M27").Value = -244.33334
$ws.Range("N27").Value = -514

# Row 35: No Risk, No Reward | Toadskin Cesti
$ws.Range("H35").Value = 1477
$ws.Range("I35").Value = 1477
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1477
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1141
$ws.Range("N35").Value = $null

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 3062.8572
$ws.Range("I61").Value = 3466.1667
$ws.Range("J61").Value = 2760.375
$ws.Range("K61").Value = 3466.1667
$ws.Range("L61").Value = 2760.375
$ws.Range("M61").Value = -3264.1667
$ws.Range("N61").Value = -3164.375

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 4000
$ws.Range("I93").Value = 3000
$ws.Range("J93").Value = 4666.6665
$ws.Range("K93").Value = 3000
$ws.Range("L93").Value = 4666.6665
$ws.Range("M93").Value = -1752
$ws.Range("N93").Value = -7162.6665

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 3062.8572
$ws.Range("I113").Value = 3466.1667
$ws.Range("J113").Value = 2760.375
$ws.Range("K113").Value = 3466.1667
$ws.Range("L113").Value = 2760.375
$ws.Range("M113").Value = -1296.1667
$ws.Range("N113").Value = -7100.375

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 3084.5715
$ws.Range("I126").Value = 2116.125
$ws.Range("J126").Value = 4375.8335
$ws.Range("K126").Value = 6348.375
$ws.Range("L126").Value = 13127.5005
$ws.Range("M126").Value = -3878.375
$ws.Range("N126").Value = -18067.5005

# Row 127: Loyal Turncoat | Saigaskin Coat of Fending
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 2845.3333
$ws.Range("I132").Value = 2045.6666
$ws.Range("J132").Value = 4444.6665
$ws.Range("K132").Value = 6136.9998
$ws.Range("L132").Value = 13333.9995
$ws.Range("M132").Value = -3606.9998
$ws.Range("N132").Value = -18393.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 123: Helping Handwear | Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 47258.92
$ws.Range("J123").Value = 47258.92
$ws.Range("L123").Value = 47258.92
$ws.Range("N123").Value = -57058.92
